$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11 formula: =A1 -> =A1*2
$ws.Range("B11").Formula = "=A1*2"

# New cell D11 = 1 (plain numeric value, default style)
$ws.Range("D11").Value = 1

# New cell E11 = "MHz" (same string as C11)
$ws.Range("E11").Value = $ws.Range("C11").Value2

# Update B12 formula: =(B11*2/(B14/1000)-16)/2 -> =(B11/(B14/1000)-16)/2
$ws.Range("B12").Formula = "=(B11/(B14/1000)-16)/2"

# New cell D12 = 10, with style matching B14 (numFmt "0" int + yellow fill)
$ws.Range("B14").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = 10

# New cell D13 = formula =B13
$ws.Range("D13").Formula = "=B13"

# New cell D14 = formula, style like B21 (numFmtId=1 only, no fill)
$ws.Range("B21").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Formula = "=D11*1000/(16+(2*D13*D12))"

# New cell E14 = "kHz" (same string as C14)
$ws.Range("E14").Value = $ws.Range("C14").Value2

$excel.CutCopyMode = 0

# Update selection to D11
$ws.Range("D11").Select()
